$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine used range extent
$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -is [string] -and $val.Contains(",")) {
            $cell.Value = $val.Replace(",", " ")
        }
    }
}
